$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    0, 2, 0, 7, 5, 4, 5, 6, 4, 4,
    4, 4, 3, 6, 2, 6, 3, 4, 4, 2,
    4, 4, 7, 7, 6, 5, 5, 9, 2, 6,
    7, 2, 4, 7, 2, 3
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
